$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1359447004608295
$ws.Range("C2").Value = 0.6566820276497696
$ws.Range("J2").Value = 0.01612903225806452
$ws.Range("P2").Value = 0.1105990783410138
$ws.Range("S2").Value = 0.08064516129032258
$ws.Range("C3").Value = 0.04013377926421405
$ws.Range("J3").Value = 0.003344481605351171
$ws.Range("P3").Value = 0.7826086956521739
$ws.Range("J4").Value = 0.06172839506172839
$ws.Range("O4").Value = 0.01234567901234568
$ws.Range("P4").Value = 0.691358024691358
$ws.Range("S4").Value = 0.2345679012345679
$ws.Range("B6").Value = 0.05797101449275362
$ws.Range("D6").Value = 0.01207729468599034
$ws.Range("E6").Value = 0.002415458937198068
$ws.Range("F6").Value = 0.04589371980676329
$ws.Range("J6").Value = 0.251207729468599
$ws.Range("O6").Value = 0.02898550724637681
$ws.Range("Q6").Value = 0.1618357487922705
$ws.Range("R6").Value = 0.05314009661835749
$ws.Range("S6").Value = 0.3864734299516908
$ws.Range("B7").Value = 0.108695652173913
$ws.Range("D7").Value = 0.02173913043478261
$ws.Range("E7").Value = 0.002415458937198068
$ws.Range("F7").Value = 0.06521739130434782
$ws.Range("J7").Value = 0.108695652173913
$ws.Range("O7").Value = 0.02657004830917874
$ws.Range("Q7").Value = 0.2028985507246377
$ws.Range("R7").Value = 0.06763285024154589
$ws.Range("S7").Value = 0.3961352657004831
$ws.Range("B8").Value = 0.093935790725327
$ws.Range("D8").Value = 0.0202140309155767
$ws.Range("E8").Value = 0.002378121284185493
$ws.Range("F8").Value = 0.07847800237812129
$ws.Range("J8").Value = 0.06302021403091558
$ws.Range("O8").Value = 0.02140309155766944
$ws.Range("Q8").Value = 0.1854934601664685
$ws.Range("R8").Value = 0.06658739595719382
$ws.Range("S8").Value = 0.4684898929845422
$ws.Range("B9").Value = 0.07021791767554479
$ws.Range("D9").Value = 0.01937046004842615
$ws.Range("F9").Value = 0.0387409200968523
$ws.Range("J9").Value = 0.0774818401937046
$ws.Range("O9").Value = 0.01937046004842615
$ws.Range("Q9").Value = 0.1864406779661017
$ws.Range("R9").Value = 0.1016949152542373
$ws.Range("S9").Value = 0.486682808716707
$ws.Range("B10").Value = 0.09703504043126684
$ws.Range("D10").Value = 0.02048517520215633
$ws.Range("E10").Value = 0.0005390835579514825
$ws.Range("F10").Value = 0.05983827493261455
$ws.Range("J10").Value = 0.1024258760107817
$ws.Range("O10").Value = 0.02210242587601078
$ws.Range("Q10").Value = 0.2280323450134771
$ws.Range("R10").Value = 0.07277628032345014
$ws.Range("S10").Value = 0.3967654986522911
$ws.Range("F11").Value = 0.001579778830963665
$ws.Range("G11").Value = 0.1469194312796208
$ws.Range("J11").Value = 0.06161137440758294
$ws.Range("K11").Value = 0.2164296998420221
$ws.Range("L11").Value = 0.5102685624012638
$ws.Range("S11").Value = 0.0631911532385466
$ws.Range("G12").Value = 0.7456647398843931
$ws.Range("J12").Value = 0.0838150289017341
$ws.Range("K12").Value = 0.0115606936416185
$ws.Range("L12").Value = 0.04046242774566474
$ws.Range("S12").Value = 0.1184971098265896
$ws.Range("G13").Value = 0.7261904761904762
$ws.Range("J13").Value = 0.130952380952381
$ws.Range("S13").Value = 0.1428571428571428
$ws.Range("F15").Value = 0.04740406320541761
$ws.Range("H15").Value = 0.1128668171557562
$ws.Range("I15").Value = 0.07900677200902935
$ws.Range("J15").Value = 0.2618510158013544
$ws.Range("K15").Value = 0.0654627539503386
$ws.Range("M15").Value = 0.01128668171557562
$ws.Range("N15").Value = 0.004514672686230248
$ws.Range("O15").Value = 0.07674943566591422
$ws.Range("S15").Value = 0.3408577878103837
$ws.Range("F16").Value = 0.02513966480446927
$ws.Range("H16").Value = 0.1675977653631285
$ws.Range("I16").Value = 0.08100558659217877
$ws.Range("J16").Value = 0.329608938547486
$ws.Range("K16").Value = 0.08659217877094973
$ws.Range("M16").Value = 0.02793296089385475
$ws.Range("N16").Value = 0.002793296089385475
$ws.Range("O16").Value = 0.06145251396648044
$ws.Range("S16").Value = 0.217877094972067
$ws.Range("F17").Value = 0.02164502164502164
$ws.Range("H17").Value = 0.1720779220779221
$ws.Range("I17").Value = 0.1038961038961039
$ws.Range("J17").Value = 0.3441558441558442
$ws.Range("K17").Value = 0.1060606060606061
$ws.Range("M17").Value = 0.01623376623376623
$ws.Range("N17").Value = 0.001082251082251082
$ws.Range("O17").Value = 0.05627705627705628
$ws.Range("S17").Value = 0.1785714285714286
$ws.Range("F18").Value = 0.01257861635220126
$ws.Range("H18").Value = 0.1635220125786163
$ws.Range("I18").Value = 0.1163522012578616
$ws.Range("J18").Value = 0.3238993710691824
$ws.Range("K18").Value = 0.1226415094339623
$ws.Range("M18").Value = 0.0220125786163522
$ws.Range("O18").Value = 0.06289308176100629
$ws.Range("S18").Value = 0.1761006289308176
$ws.Range("F19").Value = 0.02460202604920405
$ws.Range("H19").Value = 0.1928364688856729
$ws.Range("I19").Value = 0.07850940665701882
$ws.Range("J19").Value = 0.2583212735166425
$ws.Range("K19").Value = 0.09985528219971057
$ws.Range("M19").Value = 0.01917510853835022
$ws.Range("N19").Value = 0.001085383502170767
$ws.Range("O19").Value = 0.06693198263386396
$ws.Range("S19").Value = 0.2586830680173661
